# Added support for MQTT messages from MQTT broker, changed main screen
$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# Typography sheet, row 4 (verdana.ttf / Default typography definition):
# Widget Wildcard Characters (G4) gets extended to include a long set of
# punctuation / symbol characters used for testing.
$wsTypography.Range("G4").Value = "!@#$%^&*()_+-=,.<>;':""[]{}\|``~"

# Wildcard Ranges (I4) gets extended to cover lowercase and uppercase letters
# in addition to digits.
$wsTypography.Range("I4").Value = "0-9,a-z,A-Z"

# Translation sheet: add a new translated text entry (row 22) for the new
# MQTT messages text box shown on the main screen.
$wsTranslation.Range("B22").Value = "SingleUseId29"
$wsTranslation.Range("C22").Value = "Default"
$wsTranslation.Range("D22").Value = "Center"
$wsTranslation.Range("E22").Value = "LTR"
$wsTranslation.Range("F22").Value = "<value>"
